$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 19 (pushes everything from the old row19 down) ---
# First, insert a blank row above the current row 19 ("Ctrl + C" row) so the
# table keeps its blank-row separators in the right places, shifting
# rows 19-28 down to 20-29.
$ws.Rows.Item(19).Insert()

# Row 18: "Ctrl + S" / "Save" -> "  + Shift" / "Swap hands (for notes at current position)"
$ws.Range("A18").Value = "'  + Shift"
$ws.Range("B18").Value = "Swap hands (for notes at current position)"

# Row 19 (new): "  + Ctrl" / "Save"
$ws.Range("A19").Value = "'  + Ctrl"
$ws.Range("B19").Value = "Save"

# The old row 19 ("blank separator") is now row 20 - keep it blank (already is).

# --- Insert a new row before what is now row 28 (the "Esc"/"Menu" row) ---
# so a new row 28 "  + Shift" / "Play From Beginning" can be added, and the
# "Esc"/"Menu" row becomes row 29.
$ws.Rows.Item(28).Insert()

$ws.Range("A28").Value = "'  + Shift"
$ws.Range("B28").Value = "Play From Beginning"

# --- Update sheet view ---
$ws.Range("A29").Select()
$excel.ActiveWindow.ScrollRow = 4

$wb.Save()
